# Append 45 new master-data rows (102-146) to the
# "master-reg_center_machine_devic" sheet, matching the pattern of the
# existing rows, then update the page setup / selection to reflect the
# state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id / machine_id pairs that repeat every 9 rows
$aVals = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$bVals = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)

$startRow = 102
$rowCount = 45
$deviceId = 3000121

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $pairIdx = $i % 9

    $arr = New-Object 'object[,]' 1,8
    $arr[0,0] = $aVals[$pairIdx]      # regcntr_id
    $arr[0,1] = $bVals[$pairIdx]      # machine_id
    $arr[0,2] = $deviceId             # device_id
    $arr[0,3] = "eng"                 # lang_code
    $arr[0,4] = $true                 # is_active
    $arr[0,5] = "superadmin()"        # cr_by
    $arr[0,6] = "now()"               # cr_dtimes
    $arr[0,7] = "now()"               # eff_dtimes

    $rangeAddr = "A" + $row + ":H" + $row
    $ws.Range($rangeAddr).Value = $arr

    $deviceId = $deviceId + 1
}

# The sheet was left in portrait page-setup orientation.
$ws.PageSetup.Orientation = 1

# Final selection: the empty row right after the new data, extended to the
# bottom/right of the sheet (mirrors Excel's "select to end" behaviour).
$lastRow = $startRow + $rowCount
$selAddr = "A" + $lastRow + ":XFD1048576"
$ws.Range($selAddr).Select() | Out-Null
